$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values are stored as text in the source data (e.g. thousands-dot
# separated prices like "59.428.79", or plain decimals like "142.44" that would
# otherwise be auto-detected as numbers). Force text format so COM keeps them as
# literal strings, matching the workbook's original inline-string cell type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.428.79"
$ws.Range("E2").Value = "  -2.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.583.72"
$ws.Range("E3").Value = "  -2.40%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("E5").Value = "  -1.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.44"
$ws.Range("E6").Value = "  -2.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  -2.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.592.30"
$ws.Range("E9").Value = "  -2.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.164"
$ws.Range("E12").Value = "  +12.40%  "

$ws.Range("E13").Value = "  +3.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.036.15"
$ws.Range("E14").Value = "  -2.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.416.13"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.09"
$ws.Range("E16").Value = "  +5.36%  "

$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.585.97"
$ws.Range("E18").Value = "  -2.53%  "

$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.07"
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.37"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.30"
$ws.Range("E24").Value = "  -5.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.471"
$ws.Range("E25").Value = "  +6.78%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0778"
$ws.Range("E29").Value = "  -2.60%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.39"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.10"
$ws.Range("E34").Value = "  -0.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.06"
$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.17"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.909"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.37"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.861"
$ws.Range("E39").Value = "  -5.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "293.02"
$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.40"
$ws.Range("E43").Value = "  +5.57%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.594"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.65"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("E48").Value = "  -2.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.82"
$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.958.00"
$ws.Range("E51").Value = "  -0.20%  "
